$wb = $excel.ActiveWorkbook

# --- table2_COPR_s_g_ls_inc_raw (sheet1): n_extinct / n_tot rows updated for EH0/EH40 ---
$wsIncRaw = $wb.Worksheets.Item("table2_COPR_s_g_ls_inc_raw")
$wsIncRaw.Range("B6").Value = 6
$wsIncRaw.Range("D6").Value = 4
$wsIncRaw.Range("B7").Value = 16
$wsIncRaw.Range("D7").Value = 16

# --- table2_COPR_s_g_ls_inc_divtot (sheet2) ---
$wsIncDivtot = $wb.Worksheets.Item("table2_COPR_s_g_ls_inc_divtot")
$wsIncDivtot.Range("B2").Value = 0.562
$wsIncDivtot.Range("D2").Value = 0.625
$wsIncDivtot.Range("B5").Value = 0.062
$wsIncDivtot.Range("D5").Value = 0.125
$wsIncDivtot.Range("B6").Value = 0.375
$wsIncDivtot.Range("D6").Value = 0.25
$wsIncDivtot.Range("B7").Value = 16
$wsIncDivtot.Range("D7").Value = 16

# --- table2_COPR_s_g_ls_inc_divext (sheet3) ---
$wsIncDivext = $wb.Worksheets.Item("table2_COPR_s_g_ls_inc_divext")
$wsIncDivext.Range("B6").Value = 6
$wsIncDivext.Range("D6").Value = 4
$wsIncDivext.Range("B7").Value = 16
$wsIncDivext.Range("D7").Value = 16

# --- table2_COPR_s_g_ls_dec_raw (sheet4) ---
$wsDecRaw = $wb.Worksheets.Item("table2_COPR_s_g_ls_dec_raw")
$wsDecRaw.Range("B6").Value = 6
$wsDecRaw.Range("D6").Value = 4
$wsDecRaw.Range("B7").Value = 16
$wsDecRaw.Range("D7").Value = 16

# --- table2_COPR_s_g_ls_dec_divtot (sheet5) ---
$wsDecDivtot = $wb.Worksheets.Item("table2_COPR_s_g_ls_dec_divtot")
$wsDecDivtot.Range("D4").Value = 0.062
$wsDecDivtot.Range("B5").Value = 0.625
$wsDecDivtot.Range("D5").Value = 0.688
$wsDecDivtot.Range("B6").Value = 0.375
$wsDecDivtot.Range("D6").Value = 0.25
$wsDecDivtot.Range("B7").Value = 16
$wsDecDivtot.Range("D7").Value = 16

# --- table2_COPR_s_g_ls_dec_divext (sheet6) ---
$wsDecDivext = $wb.Worksheets.Item("table2_COPR_s_g_ls_dec_divext")
$wsDecDivext.Range("B6").Value = 6
$wsDecDivext.Range("D6").Value = 4
$wsDecDivext.Range("B7").Value = 16
$wsDecDivext.Range("D7").Value = 16

# Final UI state: the inc_raw sheet is the active tab with B3 selected
$wsIncRaw.Activate()
$wsIncRaw.Range("B3").Select()
